# Iceland Premier League workbook update
# Swaps the full data (columns B:AD) between several pairs of rows.
# Column A (sequential row index) is left untouched on each row.
#
# The pairs below correspond to matches whose records were
# transposed between rows during the 28-05-2024 base update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(5, 7),
    @(14, 16),
    @(33, 34),
    @(45, 46),
    @(54, 55),
    @(62, 63),
    @(145, 147)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
